$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "TC007_BusquedaDeTrenes_SoloIda"
$ws.Range("B8").Value = "Barcelona"
$ws.Range("C8").Value = "Sevilla"
$ws.Range("D8").Value = "//button[text()='24']"

$ws.Range("A9").Value = "TC008_BusquedaDeTrenes_IdaYVuelta"
$ws.Range("B9").Value = "Ciudad Real"
$ws.Range("C9").Value = "Córdoba"
$ws.Range("D9").Value = "//button[text()='19']"
$ws.Range("E9").Value = "//button[text()='24']"

$ws.Range("A10").Value = "TC0009_BusquedaDeTrenes_IdaYVuelta_MasRapido_IdaYVueltaAlMismoTiempo"
$ws.Range("B10").Value = "Alicante"
$ws.Range("C10").Value = "Madrid"
$ws.Range("D10").Value = "//button[text()='17']"
$ws.Range("E10").Value = "//button[text()='19']"

$ws.Range("A11").Value = "TC0010_BusquedaDeTrenes_IdaYVuelta_MasRapido_ReservarAsistenciaEspecial_ModalidadReducida"
$ws.Range("B11").Value = "Alicante"
$ws.Range("C11").Value = "Madrid"
$ws.Range("D11").Value = "//button[text()='17']"
$ws.Range("E11").Value = "//button[text()='19']"
$ws.Range("F11").Value = "Dylan"
$ws.Range("G11").Value = "Huarcaya"
$ws.Range("H11").Value = 16
$ws.Range("I11").Value = "//span[text()='enero']"
$ws.Range("J11").Value = 2003
$ws.Range("K11").Value = "65004204V"

$ws.Range("K11").Select()
$excel.ActiveWindow.Zoom = 145
